$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column O (28-jun) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell O1 -> copy the style of N1 (bold/centered header with border)
$wsPrix.Range("O1").Value = "28-jun"
$wsPrix.Range("N1").Copy() | Out-Null
$wsPrix.Range("O1").PasteSpecial(-4122) | Out-Null

$oValues = @{
    2  = 95.11
    3  = 78.5
    4  = 65.56999999999999
    5  = 35.07
    6  = 51.29
    7  = 44.3
    8  = 31.81
    9  = 47.38
    10 = 27.8
    11 = 2.34
    12 = 0.01
    13 = 0
    14 = -0.01
    15 = -0.02
    16 = -0.02
    17 = -0.02
    18 = 2.12
    19 = 10.26
    20 = 65
    21 = 96.26000000000001
    22 = 110
    23 = 105.09
    24 = 117.48
    25 = 103.5
}

foreach ($row in $oValues.Keys) {
    $wsPrix.Cells.Item($row, 15).Value = $oValues[$row]
}

# --- Sheet "Gaz": add row 12 (2025-06-26, 32.625) ---
$wsGaz = $wb.Worksheets.Item("Gaz")

# Set as text first (prevents auto date conversion), then re-apply the
# plain (unstyled) format used by the other date cells in column A.
$wsGaz.Range("A12").NumberFormat = "@"
$wsGaz.Range("A12").Value = "2025-06-26"
$wsGaz.Range("A11").Copy() | Out-Null
$wsGaz.Range("A12").PasteSpecial(-4122) | Out-Null
$wsGaz.Range("B12").Value = 32.625

# --- Sheet "CO2": add row 12 (2025-06-26, 69.45999999999999) ---
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A12").NumberFormat = "@"
$wsCO2.Range("A12").Value = "2025-06-26"
$wsCO2.Range("A11").Copy() | Out-Null
$wsCO2.Range("A12").PasteSpecial(-4122) | Out-Null
$wsCO2.Range("B12").Value = 69.45999999999999

Write-Host "done"
